$d = $word.ActiveDocument

function Find-ParagraphByText($pattern) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $pattern) {
            return $p
        }
    }
    return $null
}

# The three paragraphs that follow "LOQ4031: Química Geral I (Requisito fraco)"
# -- a blank paragraph, the "Ver no Jupiter..." line and the "© 2020 ..." footer
# line -- are removed, leaving only the single blank paragraph that originally
# sat before the trailing page-break paragraph.

# 1) Strip the text runs of the two footer paragraphs, leaving them blank.
$jupiter = Find-ParagraphByText("*Ver no Jupiter Salvar em pdf Salvar em docx*")
if ($jupiter -ne $null) {
    $jupiter.Range.Delete()
}

$copyright = Find-ParagraphByText("*Contact: luizeleno@usp.br*")
if ($copyright -ne $null) {
    $copyright.Range.Delete()
}

# 2) Re-locate the requirement paragraph and collapse the now-blank
#    paragraph marks that follow it, keeping exactly one blank paragraph.
$requisito = Find-ParagraphByText("LOQ4031: Química Geral I*")
$blank1 = $requisito.Next()
$blank2 = $blank1.Next()
$blank3 = $blank2.Next()

$blank3.Range.Delete()
$blank2.Range.Delete()
